# Dallas roster update: rows 8/9 (JaVale McGee <-> Frank Ntilikina) and
# rows 17/18 (Markieff Morris <-> Justin Holiday) swap their entire player
# records (columns B..K), while the "No." column (A) stays put on its row.
#
# Column I ("Exp") holds values that are sometimes purely numeric-looking
# text (e.g. "14", "5", "9", "11") mixed with non-numeric text ("R") in the
# same column, so the source file stores them as text. Plain `.Value =`
# assignment of a numeric-looking string auto-converts to a real Number
# (like typing it into Excel would), so for that column we briefly force a
# text number format, assign, then restore the cell's style so no visible
# formatting residue is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $text) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

function Swap-Range($ws, $ref1, $ref2) {
    $c1 = $ws.Range($ref1)
    $c2 = $ws.Range($ref2)
    $v1 = $c1.Value()
    $v2 = $c2.Value()
    $c1.Value = $v2
    $c2.Value = $v1
}

# --- Rows 8 <-> 9 (No. column A8/A9 untouched) ---
Swap-Range $ws "B8" "B9"
Swap-Range $ws "C8" "C9"
Swap-Range $ws "D8" "D9"
Swap-Range $ws "E8" "E9"
Swap-Range $ws "F8" "F9"
Swap-Range $ws "G8" "G9"
Swap-Range $ws "H8" "H9"
Swap-Range $ws "K8" "K9"

# Column I ("Exp") needs to stay text-typed even though the values look
# numeric.
Set-TextValue $ws "I8" "5"
Set-TextValue $ws "I9" "14"

# College: Frank Ntilikina (now row 8) has no college -> blank; JaVale
# McGee (now row 9) attended Nevada.
$ws.Range("J8").Value = ""
$ws.Range("J9").Value = "Nevada"

# --- Rows 17 <-> 18 (No. column A17/A18 untouched) ---
Swap-Range $ws "C17" "C18"
Swap-Range $ws "D17" "D18"
Swap-Range $ws "E17" "E18"
Swap-Range $ws "F17" "F18"
Swap-Range $ws "G17" "G18"
Swap-Range $ws "H17" "H18"
Swap-Range $ws "J17" "J18"
Swap-Range $ws "K17" "K18"

# Column I for rows 17/18 also needs text-forcing.
Set-TextValue $ws "I17" "9"
Set-TextValue $ws "I18" "11"

# Jersey numbers for rows 17/18 were previously blank; now populated.
$ws.Range("B17").Value = 0
$ws.Range("B18").Value = 13
